$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new row of data (row 28)
$ws.Range("B28").Value = "User password and IsActive flag, Form permissions"
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 40865

# Copy the date number format from the row above (reuses existing style
# instead of creating a new custom number format).
$ws.Range("D27").Copy()
$ws.Range("D28").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update the SUM formula range stays the same (SUM(C4:C564)) - it already
# covers the new row, but make sure the formula/value is recalculated.
$ws.Range("C3").Formula = "=SUM(C4:C564)"

# Update the selected cell to C28 as in the target workbook
$ws.Range("C28").Select()

$wb.Application.Calculate()
